$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of workout/nutrition log data (rows 5-12)
$data = @(
    @{ Row=5;  A=45231; B="Cycling"; C=2.5; D=40; E=220; F="Tempo";      G=$null; H=$null; I=$null; J=$null;          K=2500; L=150; M=300; N=80; O=50; P=2000; Q=8;   R=70;   S=50; T=8;  U=9;  V=2; W=30; X=60; Y=250; Z=8 }
    @{ Row=6;  A=45232; B="Run";     C=$null; D=$null; E=$null; F=$null; G=1;     H=8;     I=6;     J="Easy";         K=2400; L=140; M=280; N=75; O=45; P=1900; Q=7.5; R=69.8; S=52; T=7;  U=8;  V=3; W=25; X=50; Y=250; Z=9 }
    @{ Row=7;  A=45233; B="Cycling"; C=1.5; D=25; E=180; F="Recovery";   G=$null; H=$null; I=$null; J=$null;          K=2300; L=130; M=260; N=70; O=40; P=1800; Q=8.5; R=69.5; S=48; T=9;  U=9;  V=1; W=20; X=40; Y=250; Z=7 }
    @{ Row=8;  A=45234; B="Run";     C=$null; D=$null; E=$null; F=$null; G=1.5;   H=12;    I=7;     J="Tempo";        K=2600; L=160; M=320; N=85; O=55; P=2100; Q=7;   R=69.2; S=55; T=6;  U=7;  V=4; W=35; X=70; Y=250; Z=10 }
    @{ Row=9;  A=45235; B="Cycling"; C=3;   D=50; E=240; F="Interval";   G=$null; H=$null; I=$null; J=$null;          K=2700; L=170; M=340; N=90; O=60; P=2200; Q=8;   R=69;   S=53; T=8;  U=8;  V=2; W=40; X=80; Y=250; Z=9 }
    @{ Row=10; A=45236; B="Run";     C=$null; D=$null; E=$null; F=$null; G=2;     H=16;    I=8;     J="Long Run";     K=2500; L=150; M=300; N=80; O=50; P=2000; Q=9;   R=68.8; S=50; T=9;  U=9;  V=1; W=30; X=60; Y=250; Z=8.5 }
    @{ Row=11; A=45237; B="Rest";    C=$null; D=$null; E=$null; F=$null; G=$null; H=$null; I=$null; J=$null;          K=2200; L=120; M=240; N=65; O=35; P=1700; Q=10;  R=68.5; S=45; T=10; U=10; V=0; W=15; X=30; Y=250; Z=6 }
    @{ Row=12; A=45238; B="Cycling"; C=2;   D=35; E=200; F="Sweet Spot"; G=$null; H=$null; I=$null; J=$null;          K=2400; L=140; M=280; N=75; O=45; P=1900; Q=7.5; R=68.3; S=51; T=7;  U=8;  V=3; W=25; X=50; Y=250; Z=8 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

foreach ($r in $data) {
    $rowNum = $r.Row
    foreach ($col in $cols) {
        $val = $r[$col]
        $cell = $ws.Range("$col$rowNum")
        if ($val -eq $null) {
            $cell.Value = ""
        } else {
            $cell.Value = $val
        }
    }
    # Date column formatting to match column A's existing date style
    $ws.Range("A$rowNum").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
